# Taking changes to local
# The "Runmode" column (C) for rows 3-7 on the "Test Suite" sheet changes
# from "N" to "Y" (so every data row now reads "Y"), which makes the
# shared-string "N" unused and it drops out of the shared strings table
# automatically on save. The active selection on that sheet also moves
# from C3 to the C2:C7 range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")
$ws.Activate()

# Update the Runmode values for rows 3 through 7 from "N" to "Y".
$ws.Range("C3:C7").Value = "Y"

# Update the selected range/active cell to match the edited state.
[void]$ws.Range("C2:C7").Select()
